$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 82
$ws.Range("F8").Value = 2118
$ws.Range("F14").Value = 4028
$ws.Range("F17").Value = 3146
$ws.Range("F18").Value = 857
$ws.Range("F21").Value = 179
$ws.Range("F22").Value = 2093
$ws.Range("F23").Value = 1196
$ws.Range("F25").Value = 1976
$ws.Range("F27").Value = 218
$ws.Range("F28").Value = 27
$ws.Range("F29").Value = 8773
$ws.Range("F30").Value = 5797
$ws.Range("F34").Value = 14
$ws.Range("F35").Value = 780
$ws.Range("F36").Value = 3479
$ws.Range("F39").Value = 404
$ws.Range("F40").Value = 53
$ws.Range("F42").Value = 173
$ws.Range("F43").Value = 4666
$ws.Range("F45").Value = 890
$ws.Range("F46").Value = 86
$ws.Range("F47").Value = 415

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 89
$ws.Range("F17").Value = 3417

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8408
$ws.Range("F3").Value = 375
$ws.Range("F4").Value = 1342

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 8408
$ws.Range("F4").Value = 375
$ws.Range("F5").Value = 1342
$ws.Range("F6").Value = 89
$ws.Range("F8").Value = 82
$ws.Range("F11").Value = 2118
$ws.Range("F15").Value = 4028
$ws.Range("F17").Value = 3146
$ws.Range("F18").Value = 857
$ws.Range("F21").Value = 179
$ws.Range("F22").Value = 2093
$ws.Range("F26").Value = 1196
$ws.Range("F28").Value = 1976
$ws.Range("F31").Value = 218
$ws.Range("F32").Value = 27
$ws.Range("F33").Value = 8773
$ws.Range("F36").Value = 780
$ws.Range("F38").Value = 404
$ws.Range("F39").Value = 53
$ws.Range("F42").Value = 173
$ws.Range("F43").Value = 890
$ws.Range("F44").Value = 86
$ws.Range("F45").Value = 415

